$p = $ppt.ActivePresentation

# --- 1. Update the cached "datetimeFigureOut" date placeholder text
#        (slide master + every slide layout) from 29/10/2018 to 1/11/2018.
$slideMaster = $p.SlideMaster

$masterShapes = $slideMaster.Shapes
for ($j = 1; $j -le $masterShapes.Count; $j++) {
    $sh = $masterShapes.Item($j)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = "1/11/2018"
    }
}

$layouts = $slideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $lytShapes = $layouts.Item($i).Shapes
    for ($j = 1; $j -le $lytShapes.Count; $j++) {
        $sh = $lytShapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "1/11/2018"
        }
    }
}

# --- 2. Slide 1 edits -------------------------------------------------
$s = $p.Slides.Item(1)

# "currentBatch" textbox -> renamed/repositioned "currBatchPointer"
$batchLabel = $s.Shapes.Item(3)
$batchLabel.Left = -28.16165542602539
$batchLabel.Top = 73.97008514404297
$batchLabel.Width = 176.99945068359375
$batchLabel.Height = 29.081260681152344
$batchLabel.TextFrame.TextRange.Text = "currBatchPointer"

# Arrow connector pointing at it -> un-flip + reposition/resize
$batchArrow = $s.Shapes.Item(5)
$batchArrow.VerticalFlip = 0
$batchArrow.Left = 58.04653549194336
$batchArrow.Top = 103.1292953491211
$batchArrow.Width = 0
$batchArrow.Height = 60.37763977050781

# Remove the old "select 1" textbox + its arrow connector (last two shapes)
$s.Shapes.Item(28).Delete()
$s.Shapes.Item(27).Delete()
